# core dog, kick core watdog when reflashing after system down.
# Adds two more Core Watchdog Timer expire-value rows (TCR[WPEXT]||TCR[WP]
# counts 40 and 39) to the e200z3 table on sheet1, extending Table2 from
# B2:E11 to B2:E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Grow the table by two rows (keeps the ListObject/table XML in sync).
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Column B - TCR[WPEXT]||TCR[WP] count values (stored as text, like the
# existing rows).
$ws.Range("B12").Value = "'40"
$ws.Range("B13").Value = "'39"

# Column C - Expire Value (2^n)
$ws.Range("C12").Value = "2^24"
$ws.Range("C13").Value = "2^25"

# Column D - Real Time(us) (2^n/80)
$ws.Range("D12").Value = "2^24/80"
$ws.Range("D13").Value = "2^25/80"

# Column E - Real Time(ms), also stored as text.
$ws.Range("E12").Value = "'209.7152"
$ws.Range("E13").Value = "'419.4304"

# Match the author's final selection/cursor position.
$ws.Range("E13").Select() | Out-Null
